# Update the second row of per-episode metric values (columns C..L) on Sheet1
# with the freshly computed results from the first half of the first
# parameter-search function (see commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1567115.623891188
$ws.Range("D2").Value = 1567138.46577982
$ws.Range("E2").Value = 1567030.788621371
$ws.Range("F2").Value = 1566926.162664846
$ws.Range("G2").Value = 1567207.703912765
$ws.Range("H2").Value = 1567253.086686487
$ws.Range("I2").Value = 1567028.905865692
$ws.Range("J2").Value = 1567165.431763784
$ws.Range("K2").Value = 1567337.728605261
$ws.Range("L2").Value = 1567233.742292779
